# Update the "Hits" counts and derived "Percentage" text for the sheets
# whose supporting-entity visualization now differentiates between
# non-numerical and numerical answers: Total Hits, Hits_entity, Hits_boolean.
# (Hits_numerical, Hits_date, Hits_string are unaffected.)

function Set-HitsRow($ws, $row, $hits, $pct) {
    # B column: plain numeric hit count.
    $ws.Cells.Item($row, 2).Value = $hits

    # D column: literal percentage text (not a parsed percentage number).
    # Force text interpretation via NumberFormat, assign, then clear the
    # formatting override so no stray style is left on the cell.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $pct
    $dCell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Total Hits ---
$ws = $wb.Sheets.Item("Total Hits")
Set-HitsRow $ws 2 1578 "54.19%"
Set-HitsRow $ws 3 3171 "54.45%"
Set-HitsRow $ws 4 4752 "54.40%"
Set-HitsRow $ws 5 6307 "54.15%"
Set-HitsRow $ws 6 7902 "54.27%"

# --- Hits_entity ---
$ws = $wb.Sheets.Item("Hits_entity")
Set-HitsRow $ws 2 761  "53.93%"
Set-HitsRow $ws 3 1523 "53.97%"
Set-HitsRow $ws 4 2295 "54.22%"
Set-HitsRow $ws 5 3049 "54.02%"
Set-HitsRow $ws 6 3819 "54.13%"

# --- Hits_boolean ---
$ws = $wb.Sheets.Item("Hits_boolean")
Set-HitsRow $ws 2 300  "52.36%"
Set-HitsRow $ws 3 605  "52.79%"
Set-HitsRow $ws 4 897  "52.18%"
Set-HitsRow $ws 5 1198 "52.27%"
Set-HitsRow $ws 6 1499 "52.32%"
